$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 4125.9443
$ws.Range("I74").Value = 4152.6
$ws.Range("J74").Value = 4115.6924
$ws.Range("K74").Value = 4152.6
$ws.Range("L74").Value = 4115.6924
$ws.Range("M74").Value = -3216.6
$ws.Range("N74").Value = -5987.6924
$ws.Range("H76").Value = 142463.75
$ws.Range("I76").Value = 280603.5
$ws.Range("J76").Value = 4324
$ws.Range("K76").Value = 280603.5
$ws.Range("L76").Value = 4324
$ws.Range("M76").Value = -280288.5
$ws.Range("N76").Value = -4954
$ws.Range("H77").Value = 4125.9443
$ws.Range("I77").Value = 4152.6
$ws.Range("J77").Value = 4115.6924
$ws.Range("K77").Value = 20763
$ws.Range("L77").Value = 20578.462
$ws.Range("M77").Value = -16083
$ws.Range("N77").Value = -29938.462
$ws.Range("H79").Value = 142463.75
$ws.Range("I79").Value = 280603.5
$ws.Range("J79").Value = 4324
$ws.Range("K79").Value = 280603.5
$ws.Range("L79").Value = 4324
$ws.Range("M79").Value = -279511.5
$ws.Range("N79").Value = -6508
$ws.Range("H138").Value = 4118.3647
$ws.Range("I138").Value = 1952.1786
$ws.Range("J138").Value = 5010.3237
$ws.Range("K138").Value = 5856.5358
$ws.Range("L138").Value = 15030.9711
$ws.Range("M138").Value = -716.5357999999997
$ws.Range("N138").Value = -25310.9711

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 29551.957
$ws.Range("I32").Value = 18254.334
$ws.Range("J32").Value = 37711.35
$ws.Range("K32").Value = 18254.334
$ws.Range("L32").Value = 37711.35
$ws.Range("M32").Value = -17967.334
$ws.Range("N32").Value = -38285.35
$ws.Range("H45").Value = 2353.0588
$ws.Range("I45").Value = 2900.375
$ws.Range("J45").Value = 1866.5555
$ws.Range("K45").Value = 2900.375
$ws.Range("L45").Value = 1866.5555
$ws.Range("M45").Value = -2523.375
$ws.Range("N45").Value = -2620.5555
$ws.Range("H61").Value = 1708.5333
$ws.Range("I61").Value = 1553.0714
$ws.Range("J61").Value = 3885
$ws.Range("K61").Value = 1553.0714
$ws.Range("L61").Value = 3885
$ws.Range("M61").Value = -1341.0714
$ws.Range("N61").Value = -4309
$ws.Range("H132").Value = 21101.072
$ws.Range("I132").Value = 24529.283
$ws.Range("J132").Value = 3579.111
$ws.Range("K132").Value = 73587.849
$ws.Range("L132").Value = 10737.333
$ws.Range("M132").Value = -71057.849
$ws.Range("N132").Value = -15797.333
$ws.Range("H136").Value = 1708.5333
$ws.Range("I136").Value = 1553.0714
$ws.Range("J136").Value = 3885
$ws.Range("K136").Value = 4659.2142
$ws.Range("L136").Value = 11655
$ws.Range("M136").Value = -2109.2142
$ws.Range("N136").Value = -16755

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 314024.22
$ws.Range("I134").Value = 371436.12
$ws.Range("J134").Value = 4000
$ws.Range("K134").Value = 1114308.36
$ws.Range("L134").Value = 12000
$ws.Range("M134").Value = -1111773.36
$ws.Range("N134").Value = -17070

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2720.3857
$ws.Range("I31").Value = 1369
$ws.Range("J31").Value = 3568.9302
$ws.Range("K31").Value = 1369
$ws.Range("L31").Value = 3568.9302
$ws.Range("M31").Value = -1074
$ws.Range("N31").Value = -4158.9302
$ws.Range("H34").Value = 2720.3857
$ws.Range("I34").Value = 1369
$ws.Range("J34").Value = 3568.9302
$ws.Range("K34").Value = 1369
$ws.Range("L34").Value = 3568.9302
$ws.Range("M34").Value = -1167
$ws.Range("N34").Value = -3972.9302
$ws.Range("H62").Value = 2907.8948
$ws.Range("I62").Value = 2350.8
$ws.Range("J62").Value = 3106.8572
$ws.Range("K62").Value = 2350.8
$ws.Range("L62").Value = 3106.8572
$ws.Range("M62").Value = -1726.8
$ws.Range("N62").Value = -4354.8572
$ws.Range("H65").Value = 2907.8948
$ws.Range("I65").Value = 2350.8
$ws.Range("J65").Value = 3106.8572
$ws.Range("K65").Value = 11754
$ws.Range("L65").Value = 15534.286
$ws.Range("M65").Value = -8634
$ws.Range("N65").Value = -21774.286
$ws.Range("H68").Value = 15000
$ws.Range("J68").Value = 15000
$ws.Range("L68").Value = 15000
$ws.Range("N68").Value = -16498
$ws.Range("H71").Value = 15000
$ws.Range("J71").Value = 15000
$ws.Range("L71").Value = 45000
$ws.Range("N71").Value = -52488

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 689.6094000000001
$ws.Range("I5").Value = 526.40424
$ws.Range("J5").Value = 1140.8235
$ws.Range("K5").Value = 1579.21272
$ws.Range("L5").Value = 3422.4705
$ws.Range("M5").Value = -1467.21272
$ws.Range("N5").Value = -3646.4705
$ws.Range("H68").Value = 1796.3334
$ws.Range("I68").Value = 1798
$ws.Range("J68").Value = 1796
$ws.Range("K68").Value = 5394
$ws.Range("L68").Value = 5388
$ws.Range("M68").Value = -4583
$ws.Range("N68").Value = -7010
$ws.Range("H71").Value = 1796.3334
$ws.Range("I71").Value = 1798
$ws.Range("J71").Value = 1796
$ws.Range("K71").Value = 16182
$ws.Range("L71").Value = 16164
$ws.Range("M71").Value = -12126
$ws.Range("N71").Value = -24276
$ws.Range("H132").Value = 1063.5
$ws.Range("I132").Value = 521.2941
$ws.Range("J132").Value = 1605.7059
$ws.Range("K132").Value = 4691.6469
$ws.Range("L132").Value = 14451.3531
$ws.Range("M132").Value = -2161.6469
$ws.Range("N132").Value = -19511.3531
$ws.Range("H135").Value = 689.6094000000001
$ws.Range("I135").Value = 526.40424
$ws.Range("J135").Value = 1140.8235
$ws.Range("K135").Value = 4737.63816
$ws.Range("L135").Value = 10267.4115
$ws.Range("M135").Value = -2202.63816
$ws.Range("N135").Value = -15337.4115

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 2348.1785
$ws.Range("I97").Value = 1524.4546
$ws.Range("K97").Value = 1524.4546
$ws.Range("M97").Value = -1028.4546
$ws.Range("H132").Value = 3484.1143
$ws.Range("I132").Value = 3364.9333
$ws.Range("J132").Value = 4199.2
$ws.Range("K132").Value = 10094.7999
$ws.Range("L132").Value = 12597.6
$ws.Range("M132").Value = -7564.7999
$ws.Range("N132").Value = -17657.6
